# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" on the Overview/zh-cn/de-de sheets and refreshes
# the associated handoff timestamps, then widens the now-longer Status
# columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status + latest HO xliff generate date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-08 05:17:30"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-08 05:17:25"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-08 05:17:30"

# --- Widen the Status columns now that "Ready for handoff" is longer than
#     "In Translation" (mirrors the workbook being re-saved after the edit).
$newStatusWidth = 16.3333333333333
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
